$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 ("TS002_AddAggrement") flag changes from "No" to "Yes"
$ws.Range("C3").Value = "Yes"

# Remove the stray hyperlink that targeted the doomed row 4 (D4) before
# the row shift makes it ambiguous
foreach ($h in @($ws.Hyperlinks)) {
    $addr = $h.Range.Address()
    if ($addr -eq '$D$4') {
        $h.Delete()
    }
}

# Drop the existing Flag list validation (covers C2:C4) before the row
# shift, then delete the third data row (TS002_AggrementVerfication)
$ws.Range("C2:C4").Validation.Delete()
$ws.Rows.Item(4).Delete()

# Re-apply the Flag list validation, now only spanning the remaining rows
$ws.Range("C2:C3").Validation.Add(3, 1, 1, '"Yes,No"') | Out-Null

# Move the active selection to C6, matching the saved view state
$null = $ws.Range("C6").Select()
